$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 8 new blank rows before the existing row 113 ("*ない形 變化 I" block),
# pushing the old rows 113-137 down to 121-145.
$ws.Rows("113:120").Insert()

# Row 113: new lesson header
$ws.Range("A113").Value = "*第18課文型"

# Row 114
$ws.Range("A114").Value = "ミラーさんは かんじ を よむ ことが できます"
$ws.Range("B114").Value = "ミラーさんは 漢字を 読む ことが できます"
$ws.Range("C114").Value = "米勒會讀漢字"

# Row 115
$ws.Range("A115").Value = "わたしの しゅみ は えいが を みる ことです"
$ws.Range("B115").Value = "わたしの 趣味は 映画を 見る ことです"
$ws.Range("C115").Value = "我的興趣是看電影"

# Row 116
$ws.Range("A116").Value = "ねる まえに、にっき を かきます"
$ws.Range("B116").Value = "寝る まえに、日記を 書きます"
$ws.Range("C116").Value = "我在睡前寫日記裡"

# Row 117
$ws.Range("A117").Value = "スキーが できますか"
$ws.Range("C117").Value = "你會滑雪嗎？"

# Row 118
$ws.Range("A118").Value = "…はい、できます。 でも、あまり じょうずじゃ ありません。"
$ws.Range("B118").Value = "…はい、できます。 でも、あまり 上手じゃ ありません。"
$ws.Range("C118").Value = "……我會，但我滑得不太好。"

# Row 119
$ws.Range("A119").Value = "マリアさんは パソコンを つかう ことが できますか。"
$ws.Range("B119").Value = "マリアさんは パソコン を 使う ことが できますか。"
$ws.Range("C119").Value = "瑪麗亞會用電腦嗎？"

# Row 120
$ws.Range("A120").Value = "…いいえ、できません。"
$ws.Range("C120").Value = "……不會。"

# Update the view so the newly added block is in frame / selected, mirroring
# the author's saved cursor position after the edit.
$ws.Range("A121").Select() | Out-Null
